$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old emoji -> new replacement text, keyed so we can look the replacement
# up once we know which emoji a given cell currently holds.
$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📙" = "+3"
    "📗" = "✅"
}

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value2
        foreach ($key in $map.Keys) {
            if ($key -eq $val) {
                $newVal = $map[$key]
                # Assign via a literal-text formula, then convert that
                # formula to a plain value with a values-only paste. This
                # keeps replacements such as "-3"/"+3" stored as genuine
                # text (no numeric re-interpretation, no quote-prefix
                # style added) instead of going through .Value directly.
                $cell.Formula = '="' + $newVal + '"'
                $cell.Copy()
                $cell.PasteSpecial(-4163)
                break
            }
        }
    }
}

$excel.CutCopyMode = 0
